$wb = $excel.ActiveWorkbook

# Apply the same header-row additions (columns I:M) to every sheet in the
# workbook: a "use custom unit" toggle column plus two unit/quantity pairs.
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Write the new header cells. The order below matches the order the
    # strings were first introduced into the shared-string table.
    $ws.Range("J1").Value = "单位1"
    $ws.Range("L1").Value = "单位2"
    $ws.Range("K1").Value = "单位1数量"
    $ws.Range("M1").Value = "单位2数量"
    $ws.Range("I1").Value = "使用自定义单位"

    # J1/K1/L1/M1 just need vertical-centered text (matches the rest of the
    # header row's plain cells).
    $ws.Range("J1:M1").VerticalAlignment = -4108

    # I1 additionally carries a dedicated font (11pt 等线) alongside the
    # vertical centering.
    $ws.Range("I1").VerticalAlignment = -4108
    $ws.Range("I1").Font.Family = 3

    # New column I is sized like the other description columns.
    $ws.Columns.Item(9).ColumnWidth = 14.43
}

# Sheet-specific selection state, matching what was left active after the edit.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Range("I2").Select()

$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$ws2.Range("I1:M1").Select()

$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate()
$ws3.Range("I1:M1").Select()

$ws1.Activate()
